$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A4").Value = -20.12229999999998
$ws.Range("A7").Value = -21.95090000000002
$ws.Range("B7").Value = 4.748200000000004
$ws.Range("B15").Value = 5.055299999999996
$ws.Range("A16").Value = -21.53149999999998
$ws.Range("C16").Value = -11.5799
$ws.Range("C19").Value = -12.33900000000001
$ws.Range("B21").Value = 10.4876
$ws.Range("B22").Value = 10.346
$ws.Range("B23").Value = 8.682600000000006
$ws.Range("A28").Value = -22.0141
$ws.Range("A29").Value = -21.34269999999997
$ws.Range("A32").Value = -21.16149999999999
$ws.Range("B34").Value = 8.312900000000003
$ws.Range("D34").Value = -8.043699999999999
$ws.Range("C36").Value = -12.6198
$ws.Range("A40").Value = -20.06829999999999
$ws.Range("B43").Value = 5.375100000000003
$ws.Range("D43").Value = -8.429499999999999
$ws.Range("B45").Value = 4.961300000000001
$ws.Range("C46").Value = -14.59799999999999
$ws.Range("D48").Value = -7.506599999999999
$ws.Range("B50").Value = 5.248099999999997
$ws.Range("C50").Value = -13.02929999999999
$ws.Range("B51").Value = 5.827099999999996
$ws.Range("A52").Value = -21.99829999999999
$ws.Range("A57").Value = -22.42830000000002
$ws.Range("A66").Value = -21.4396
$ws.Range("B66").Value = 5.231399999999999
$ws.Range("B67").Value = 5.213500000000002
$ws.Range("D70").Value = -7.537800000000002
$ws.Range("D73").Value = -7.696499999999994
$ws.Range("B79").Value = 10.0206
$ws.Range("B84").Value = 5.676900000000001
$ws.Range("D87").Value = -8.608999999999995
$ws.Range("B92").Value = 4.731099999999997
$ws.Range("D92").Value = -6.219700000000001
$ws.Range("C95").Value = -12.2447
$ws.Range("B97").Value = 6.589399999999997
$ws.Range("C97").Value = -12.4114
$ws.Range("A100").Value = -22.0963
$ws.Range("D101").Value = -8.032999999999996
